$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# overworld_1_hint: replace placeholder text with real hint, bump VoiceDuration 5 -> 10
$ws.Range("B110").Value = "In a region with temperate climate, the highest temperature and humidity tends to be during summer."
$ws.Range("C110").Value = 10

# overworld_2_hint
$ws.Range("B154").Value = "You’re certain to find strong winds in a tropical region during autumn in the Pacific Ocean."
$ws.Range("C154").Value = 10

# overworld_3_hint
$ws.Range("B174").Value = "Northern Africa is a good desert region for these frogs. Just make sure to pick a season that’s not too hot!"
$ws.Range("C174").Value = 10

# overworld_4_hint - also loses its wrap-text/vertical-center style in the target
$ws.Range("B186").Value = "We want a mountainous area for this one. How about in South America during summer, when it is cool?"
$ws.Range("C186").Value = 10
$ws.Range("B186").Style = "Normal"

# Update the visible window / selection to match the author's final cursor position
$ws.Application.ActiveWindow.ScrollRow = 161
[void]$ws.Range("B186").Select()
